# update field beauty - MDLWL
# Cell A2 on Sheet1 is re-generated to a fresh unique id value, cycling
# through a few freshly-generated candidate ids before landing on the
# final one (mirrors the id-generation history captured in the
# sharedStrings table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "CA-7QX1CIST"
$ws.Range("A2").Value = "CA-7REDN9KS"
$ws.Range("A2").Value = "CA-UFN14XHY"
$ws.Range("A2").Value = "CA-9RL3N43P"
$ws.Range("A2").Value = "CA-MWQYTQLX"
